$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fill in the "position" column (G) for the three new employees
#     first (matches the order the shared strings were originally authored in).
$ws.Range("G4").Value2 = "L1"
$ws.Range("G5").Value2 = "fired"
$ws.Range("G6").Value2 = "L1"

# --- Step 2: row 4 (id 3 - Papov / 1 / Adel) --------------------------------
$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value2 = "Papov"
# "1" looks numeric, so force text entry (matches the t="s" cell in the sheet)
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value2 = "1"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value2 = "Adel"
$ws.Range("E4").Value2 = 27489
$ws.Range("F4").Value2 = 1

# --- Step 3: row 5 (id 4 - сотрудник / 9 / сотрудник) -----------------------
$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value2 = "сотрудник"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value2 = "9"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value2 = "сотрудник"
$ws.Range("E5").Value2 = 31977
$ws.Range("F5").Value2 = 2

# --- Step 4: row 6 (id 5 - 9 / 9 / 9) ---------------------------------------
$ws.Range("A6").Value2 = 5
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value2 = "9"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value2 = "9"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "9"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = 31977
$ws.Range("F6").Value2 = 3
